# Feria Lagunitas de Puerto Montt - Poroto verde
# A new weekly price record is added at the top of the data block (row 7),
# pushing all subsequent records down by one row. The previous last record
# (old row 67) ends up duplicated into the new last row (68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7..67 down by inserting a new blank row at row 7.
# This naturally creates the new row 68 (copy of old row 67) and leaves
# row 7 empty, ready to receive the new record's data.
$ws.Rows.Item(7).Insert()

# Populate the new record in row 7.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C7").Value = 'Los Lagos'
$ws.Range("D7").Value = 44649
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112031
$ws.Range("G7").Value = 'Poroto verde'
$ws.Range("H7").Value = 'Magnum'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 30000
$ws.Range("N7").Value = '$/saco 25 kilos'
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 1200
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = 'Hortaliza'
